$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48)
$colA = @(
    '长沙市开福区高岭香江国际城充电站建设项目',
    '长沙市开福区高岭香江国际城充电站建设项目',
    '长沙市开福区高岭香江国际城充电站建设项目',
    '长沙市开福区高岭香江国际城充电站建设项目',
    '长沙市开福区高岭香江国际城充电站建设项目',
    '长沙市开福区高岭香江国际城充电站建设项目',
    '长沙市开福区高岭香江国际城充电站建设项目',
    '长沙市开福区高岭香江国际城充电站建设项目',
    '长沙市开福区高岭香江国际城充电站建设项目',
    '长沙市开福区高岭香江国际城充电站建设项目',
    '长沙市开福区高岭香江国际城充电站建设项目',
    '长沙市开福区高岭香江国际城充电站建设项目',
    '长沙市开福区高岭香江国际城充电站建设项目',
    '飞狐四方坪南区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪南区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪东区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪南区充电站',
    '飞狐四方坪南区充电站',
    '飞狐四方坪东区充电站',
    '飞狐四方坪南区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪南区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪东区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪东区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪东区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪东区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站',
    '飞狐四方坪西区充电站'
)
$colB = @(
    '108号直流',
    '107号直流',
    '208号直流',
    '206号直流',
    '209号直流',
    '111号直流',
    '104号直流',
    '109号直流',
    '204号直流',
    '106号直流',
    '211号直流',
    '203号直流',
    '307号直流',
    '9176699368200101',
    '9176699400501202',
    '9176699400500303',
    '9176699400500501',
    '9176699400501205',
    '9176699368200203',
    '9176699400500502',
    '9176699442100202',
    '9176699400501102',
    '9176699400501304',
    '9176699368200103',
    '9176699368200201',
    '9176699442100402',
    '9176699368200306',
    '9176699400500201',
    '9176699355900102',
    '9176699368200406',
    '9176699400500304',
    '9176699425700302',
    '9176699400501302',
    '9176699400500404',
    '9176699400500302',
    '9176699425700301',
    '9176699400500205',
    '9176699435600102',
    '9176699400500802',
    '9176699400500305',
    '9176699442100302',
    '9176699400500104',
    '9176699400501203',
    '9176699400500102',
    '9176699400500204',
    '9176699400500601',
    '9176699400500504'
)
$colC = @(
    46056.45337962963,
    46056.507268518515,
    46056.52923611111,
    46056.54145833333,
    46056.56079861111,
    46056.57104166667,
    46056.574375,
    46056.59936342593,
    46056.63670138889,
    46056.65008101852,
    46056.67521990741,
    46056.759108796294,
    46056.78325231482,
    46055.06685185185,
    46055.249976851854,
    46055.54478009259,
    46055.552835648145,
    46055.573900462965,
    46055.80465277778,
    46056.17390046296,
    46056.22856481482,
    46056.24049768518,
    46056.343194444446,
    46056.51405092593,
    46056.52486111111,
    46056.52966435185,
    46056.536631944444,
    46056.539351851854,
    46056.539560185185,
    46056.54337962963,
    46056.54809027778,
    46056.54990740741,
    46056.55813657407,
    46056.564155092594,
    46056.576585648145,
    46056.583078703705,
    46056.58331018518,
    46056.592685185184,
    46056.60065972222,
    46056.606574074074,
    46056.652395833335,
    46056.673854166664,
    46056.678819444445,
    46056.745358796295,
    46056.770462962966,
    46056.805,
    46056.82576388889
)
$colD = @(
    46057.35747685185,
    46057.35747685185,
    46057.35747685185,
    46057.35747685185,
    46057.35747685185,
    46057.35747685185,
    46057.35747685185,
    46057.35747685185,
    46057.35747685185,
    46057.35747685185,
    46057.35747685185,
    46057.35747685185,
    46057.35747685185,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375,
    46057.344375
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
}

$ws.Range("E15").Select()